$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update weighted grade (column C) values per diff
$ws.Range("C3").Value = 56.23
$ws.Range("C5").Value = 60.03
$ws.Range("C6").Value = 60.03
$ws.Range("C7").Value = 75.69
$ws.Range("C8").Value = 29.7
$ws.Range("C11").Value = 12.01
$ws.Range("C15").Value = 21.85
$ws.Range("C18").Value = 56.23
$ws.Range("C23").Value = 21.85

# Student 74311 (row 24): update weighted grade and feedback text
$ws.Range("C24").Value = 86.91
$ws.Range("H24").Value = "Excellence! Outstanding work across almost all skills. Keep up the great work!"
